{"js": "// The document is a two-digit multiplication worksheet: a 5x5 table\n// whose first row plus every other row (25 cells total) holds an\n// \"AxB=\" expression. The edit swaps each expression (in document\n// order, reading the table left-to-right / top-to-bottom) for a new\n// one, per the diff. A couple of the new values collide with other\n// cells' original values (e.g. cell 11 becomes \"25\u00d740=\", which is\n// cell 24's *original* text), so we must not do a naive\n// search-current-text/replace loop -- a later search for the\n// original \"25\u00d740=\" would match the just-inserted one instead of the\n// real target. Instead, gather every \"=\" hit up front (one per cell,\n// in document order) before making any edits, then replace each\n// cell's whole paragraph text by position by zipping against the\n// ordered list of new values.\nconst oldValues = [\n  \"85\u00d798=\", \"68\u00d748=\", \"43\u00d780=\", \"92\u00d717=\", \"84\u00d771=\",\n  \"90\u00d756=\", \"68\u00d728=\", \"95\u00d731=\", \"49\u00d762=\", \"94\u00d792=\",\n  \"14\u00d779=\", \"56\u00d778=\", \"31\u00d711=\", \"85\u00d746=\", \"11\u00d795=\",\n  \"67\u00d726=\", \"30\u00d769=\", \"21\u00d719=\", \"93\u00d734=\", \"37\u00d746=\",\n  \"99\u00d754=\", \"89\u00d797=\", \"67\u00d781=\", \"25\u00d740=\", \"50\u00d753=\",\n];\n\nconst newValues = [\n  \"41\u00d734=\", \"39\u00d767=\", \"48\u00d772=\", \"30\u00d744=\", \"90\u00d741=\",\n  \"30\u00d759=\", \"64\u00d726=\", \"16\u00d784=\", \"54\u00d767=\", \"79\u00d731=\",\n  \"25\u00d740=\", \"32\u00d798=\", \"13\u00d772=\", \"93\u00d732=\", \"42\u00d788=\",\n  \"14\u00d735=\", \"58\u00d743=\", \"84\u00d738=\", \"48\u00d756=\", \"54\u00d741=\",\n  \"52\u00d789=\", \"12\u00d774=\", \"64\u00d767=\", \"77\u00d761=\", \"50\u00d768=\",\n];\n\nconst body = context.document.body;\n\n// Every target cell ends with \"=\", and nothing else in the document\n// does, so this single search finds all 25 cells in document order.\nconst results = body.search(\"=\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== oldValues.length) {\n  throw new Error(\n    `Expected ${oldValues.length} \"=\" matches, found ${results.items.length}`\n  );\n}\n\n// Resolve each hit to its containing paragraph (the whole \"AxB=\" cell\n// text) before any writes happen, so later edits can't shift/confuse\n// earlier anchors.\nconst paragraphs = results.items.map((item) => item.paragraphs.getFirst());\nfor (const p of paragraphs) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.length; i++) {\n  const actual = paragraphs[i].text;\n  if (actual !== oldValues[i]) {\n    throw new Error(\n      `Cell ${i}: expected \"${oldValues[i]}\" but found \"${actual}\"`\n    );\n  }\n  paragraphs[i].getRange().insertText(newValues[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document is a two-digit multiplication worksheet: a single 5-\n# column table whose content rows (1, 5, 10, 15, 20 - every other row\n# across 20 total rows) each hold 5 \"AxB=\" expressions, 25 in all.\n# The edit swaps each expression for a new one, per the diff. Several\n# of the new values collide with other cells' original values (e.g.\n# row 10 / col 1 becomes \"25x40=\", which is row 20 / col 4's original\n# text), so rather than searching the document by the OLD text (which\n# would risk matching a value some earlier step just inserted), this\n# addresses every target cell directly by its fixed (row, column)\n# position in the table, exactly like the source OOXML table-cell\n# diff does.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$contentRows = @(1, 5, 10, 15, 20)\n\n$oldValues = @(\n  @(\"85\u00d798=\", \"68\u00d748=\", \"43\u00d780=\", \"92\u00d717=\", \"84\u00d771=\"),\n  @(\"90\u00d756=\", \"68\u00d728=\", \"95\u00d731=\", \"49\u00d762=\", \"94\u00d792=\"),\n  @(\"14\u00d779=\", \"56\u00d778=\", \"31\u00d711=\", \"85\u00d746=\", \"11\u00d795=\"),\n  @(\"67\u00d726=\", \"30\u00d769=\", \"21\u00d719=\", \"93\u00d734=\", \"37\u00d746=\"),\n  @(\"99\u00d754=\", \"89\u00d797=\", \"67\u00d781=\", \"25\u00d740=\", \"50\u00d753=\")\n)\n\n$newValues = @(\n  @(\"41\u00d734=\", \"39\u00d767=\", \"48\u00d772=\", \"30\u00d744=\", \"90\u00d741=\"),\n  @(\"30\u00d759=\", \"64\u00d726=\", \"16\u00d784=\", \"54\u00d767=\", \"79\u00d731=\"),\n  @(\"25\u00d740=\", \"32\u00d798=\", \"13\u00d772=\", \"93\u00d732=\", \"42\u00d788=\"),\n  @(\"14\u00d735=\", \"58\u00d743=\", \"84\u00d738=\", \"48\u00d756=\", \"54\u00d741=\"),\n  @(\"52\u00d789=\", \"12\u00d774=\", \"64\u00d767=\", \"77\u00d761=\", \"50\u00d768=\")\n)\n\nfor ($ri = 0; $ri -lt $contentRows.Count; $ri++) {\n  $row = $contentRows[$ri]\n  for ($ci = 0; $ci -lt 5; $ci++) {\n    $col = $ci + 1\n    $cell = $t.Cell($row, $col)\n    $expected = $oldValues[$ri][$ci]\n    $actual = $cell.Range.Text\n    if ($actual -ne ($expected + \"`r`a\")) {\n      throw \"Row $row Col $col`: expected '$expected' but found '$actual'\"\n    }\n    $cell.Range.Text = $newValues[$ri][$ci]\n  }\n}\n"}
